# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.860.74'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.474.91'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.68%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.05'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.15'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.473.27'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.68%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -2.78%  '
$ws.Range("E10").Value = '  -3.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.59'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000212'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.059.27'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.29'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.466.43'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.873.12'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.40'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.27'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.99'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '431.99'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.604'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -6.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.93'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.608.84'
$ws.Range("D26").ClearFormats()
$ws.Range("E27").Value = '  -8.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.76'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -7.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.36'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -8.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.48'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.78%  '
$ws.Range("E31").Value = '  -7.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.24'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.462.18'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.92'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.86'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '174.47'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0880'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.44%  '
$ws.Range("E43").Value = '  -12.25%  '
$ws.Range("E44").Value = '  -4.21%  '
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.29'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.71'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -7.59%  '
$ws.Range("E48").Value = '  -8.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.42'
$ws.Range("D49").ClearFormats()
$ws.Range("E50").Value = '  -9.79%  '
$ws.Range("E51").Value = '  -5.02%  '
